$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D="'42.136.48"; E="  -1.95%  " }
    3 = @{ D="'2.269.20"; E="  -2.90%  " }
    4 = @{ E="  -0.02%  " }
    5 = @{ D="'297.69"; E="  -3.03%  " }
    6 = @{ D="'94.70"; E="  -6.67%  " }
    7 = @{ E="  +0.03%  " }
    8 = @{ D="'0.491"; E="  -3.66%  " }
    9 = @{ E="  -3.69%  " }
    10 = @{ D="'33.33"; E="  -4.31%  " }
    12 = @{ D="'48.22"; E="  -8.13%  " }
    13 = @{ E="  -0.21%  " }
    14 = @{ D="'6.64"; E="  -2.96%  " }
    15 = @{ D="'15.71"; E="  -0.52%  " }
    16 = @{ D="'2.620.12"; E="  -3.10%  " }
    17 = @{ D="'2.266.23"; E="  -2.14%  " }
    18 = @{ D="'0.779"; E="  -5.74%  " }
    19 = @{ D="'42.074.09" }
    20 = @{ D="'11.58"; E="  -1.33%  " }
    21 = @{ D="'0.0₃0889"; E="  -2.38%  " }
    22 = @{ E="  -3.19%  " }
    23 = @{ D="'66.47"; E="  -2.63%  " }
    24 = @{ D="'233.25"; E="  -1.39%  " }
    25 = @{ D="'1.96"; E="  -2.34%  " }
    26 = @{ E="  +0.22%  " }
    27 = @{ E="  -4.21%  " }
    28 = @{ D="'23.90"; E="  -6.66%  " }
    29 = @{ D="'2.28"; E="  -1.92%  " }
    30 = @{ D="'168.05"; E="  +3.94%  " }
    31 = @{ D="'34.06"; E="  -3.41%  " }
    32 = @{ D="'9.08"; E="  -2.33%  " }
    33 = @{ E="  -0.04%  " }
    34 = @{ D="'4.89"; E="  -4.31%  " }
    35 = @{ D="'4.49"; E="  -3.43%  " }
    36 = @{ D="'2.32"; E="  -4.99%  " }
    37 = @{ B="Celestia"; C="https://coinranking.com/coin/YQcD0lBl7+celestia-tia"; D="'16.44"; E="  -5.88%  " }
    38 = @{ B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="'0.0687"; E="  -5.10%  " }
    39 = @{ E="  -4.45%  " }
    40 = @{ D="'0.0986"; E="  -3.09%  " }
    41 = @{ E="  -3.45%  " }
    42 = @{ E="  -7.05%  " }
    43 = @{ D="'2.44"; E="  -5.43%  " }
    44 = @{ D="'1.959.43"; E="  -3.18%  " }
    45 = @{ D="'0.0277"; E="  -2.51%  " }
    46 = @{ D="'17.36"; E="  -7.87%  " }
    47 = @{ D="'9.55"; E="  -6.29%  " }
    48 = @{ E="  -5.15%  " }
    49 = @{ D="'2.493.26"; E="  -2.68%  " }
    50 = @{ D="'51.90"; E="  -7.12%  " }
    51 = @{ D="'4.54"; E="  -3.71%  " }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $addr = "$col$r"
        $ws.Range($addr).Value = $row[$col]
    }
}
